$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 44.333332
$ws.Range("I39").Value = 24.5
$ws.Range("J39").Value = 57.555557
$ws.Range("K39").Value = 73.5
$ws.Range("L39").Value = 172.666671
$ws.Range("M39").Value = 222.5
$ws.Range("N39").Value = -764.666671

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7497.5
$ws.Range("I116").Value = 992
$ws.Range("K116").Value = 992
$ws.Range("M116").Value = 2450

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 68330.47
$ws.Range("J129").Value = 102257.7
$ws.Range("L129").Value = 306773.1
$ws.Range("N129").Value = -316773.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2697397.8
$ws.Range("I132").Value = 3176602.8
$ws.Range("J132").Value = 1869
$ws.Range("K132").Value = 9529808.399999999
$ws.Range("L132").Value = 5607
$ws.Range("M132").Value = -9527278.399999999
$ws.Range("N132").Value = -10667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15285.527
$ws.Range("I32").Value = 15124.807
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 15124.807
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -14837.807
$ws.Range("N32").Value = -20574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1485.5385
$ws.Range("I74").Value = 1789
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1789
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -915
$ws.Range("N74").Value = -2748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1485.5385
$ws.Range("I77").Value = 1789
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 8945
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -4577
$ws.Range("N77").Value = -13736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1799
$ws.Range("I122").Value = 1653.0416
$ws.Range("J122").Value = 2966.6667
$ws.Range("K122").Value = 4959.1248
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("M122").Value = -2509.1248
$ws.Range("N122").Value = -13800.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5481.8965
$ws.Range("I132").Value = 5796.2554
$ws.Range("J132").Value = 4138.727
$ws.Range("K132").Value = 17388.7662
$ws.Range("L132").Value = 12416.181
$ws.Range("M132").Value = -14858.7662
$ws.Range("N132").Value = -17476.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4806.8203
$ws.Range("I134").Value = 5509.857
$ws.Range("J134").Value = 3017.2727
$ws.Range("K134").Value = 16529.571
$ws.Range("L134").Value = 9051.8181
$ws.Range("M134").Value = -13994.571
$ws.Range("N134").Value = -14121.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3626771.5
$ws.Range("I31").Value = 2723.3125
$ws.Range("J31").Value = 11910310
$ws.Range("K31").Value = 2723.3125
$ws.Range("L31").Value = 11910310
$ws.Range("M31").Value = -2428.3125
$ws.Range("N31").Value = -11910900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3626771.5
$ws.Range("I34").Value = 2723.3125
$ws.Range("J34").Value = 11910310
$ws.Range("K34").Value = 2723.3125
$ws.Range("L34").Value = 11910310
$ws.Range("M34").Value = -2521.3125
$ws.Range("N34").Value = -11910714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5003378.5
$ws.Range("I132").Value = 3318.25
$ws.Range("J132").Value = 7356348
$ws.Range("K132").Value = 9954.75
$ws.Range("L132").Value = 22069044
$ws.Range("M132").Value = -7424.75
$ws.Range("N132").Value = -22074104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 140
$ws.Range("I12").Value = 126.35714
$ws.Range("J12").Value = 151.23529
$ws.Range("K12").Value = 379.07142
$ws.Range("L12").Value = 453.70587
$ws.Range("M12").Value = -206.07142
$ws.Range("N12").Value = -799.70587

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14316.987
$ws.Range("I131").Value = 5852.722
$ws.Range("J131").Value = 16735.35
$ws.Range("K131").Value = 17558.166
$ws.Range("L131").Value = 50206.05
$ws.Range("M131").Value = -12518.166
$ws.Range("N131").Value = -60286.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 44850400
$ws.Range("I137").Value = 37049604
$ws.Range("J137").Value = 48041636
$ws.Range("K137").Value = 111148812
$ws.Range("L137").Value = 144124908
$ws.Range("M137").Value = -111143712
$ws.Range("N137").Value = -144135108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 52073164
$ws.Range("I70").Value = 65656252
$ws.Range("J70").Value = 4651.5
$ws.Range("K70").Value = 65656252
$ws.Range("L70").Value = 4651.5
$ws.Range("M70").Value = -65655982
$ws.Range("N70").Value = -5191.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 52073164
$ws.Range("I73").Value = 65656252
$ws.Range("J73").Value = 4651.5
$ws.Range("K73").Value = 65656252
$ws.Range("L73").Value = 4651.5
$ws.Range("M73").Value = -65655316
$ws.Range("N73").Value = -6523.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2028.5
$ws.Range("I97").Value = 1863.5625
$ws.Range("J97").Value = 2468.3333
$ws.Range("K97").Value = 1863.5625
$ws.Range("L97").Value = 2468.3333
$ws.Range("M97").Value = -1367.5625
$ws.Range("N97").Value = -3460.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5742.148
$ws.Range("I132").Value = 6193.3477
$ws.Range("J132").Value = 3147.75
$ws.Range("K132").Value = 18580.0431
$ws.Range("L132").Value = 9443.25
$ws.Range("M132").Value = -16050.0431
$ws.Range("N132").Value = -14503.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 585.5
$ws.Range("I22").Value = 657.4
$ws.Range("J22").Value = 465.66666
$ws.Range("K22").Value = 657.4
$ws.Range("L22").Value = 465.66666
$ws.Range("M22").Value = -362.4
$ws.Range("N22").Value = -1055.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 585.5
$ws.Range("I27").Value = 657.4
$ws.Range("J27").Value = 465.66666
$ws.Range("K27").Value = 657.4
$ws.Range("L27").Value = 465.66666
$ws.Range("M27").Value = -550.4
$ws.Range("N27").Value = -679.66666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1365.9259
$ws.Range("I132").Value = 1273.289
$ws.Range("J132").Value = 1829.1111
$ws.Range("K132").Value = 3819.867
$ws.Range("L132").Value = 5487.3333
$ws.Range("M132").Value = -1289.867
$ws.Range("N132").Value = -10547.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 986.03925
$ws.Range("I136").Value = 1028.7609
$ws.Range("J136").Value = 593
$ws.Range("K136").Value = 3086.2827
$ws.Range("L136").Value = 1779
$ws.Range("M136").Value = -536.2826999999997
$ws.Range("N136").Value = -6879
